# Update: Error handling + Validation
#
# Target the "Transactions" sheet explicitly (the workbook's active sheet
# may be "Metadata" when it is loaded).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Remove the "Balance" column (E) entirely - header + data - shifting
# everything left so the used range becomes A1:D2 (then A1:D5 once the
# new rows below are added).
$ws.Range("E1:E2").Delete()

# Helper-free approach: write each new cell as literal TEXT (not a parsed
# date/number) by pre-formatting the cell as Text ("@") before assigning
# the value, then resetting the style back to "Normal" so no stray number
# format is left attached to the cell (matches the source data, which has
# no style on these rows).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 3
Set-TextValue $ws.Range("A3") "03-05-2024"
$ws.Range("B3").Value = "ATM Withdrawal"
Set-TextValue $ws.Range("C3") "5,000.00"
$ws.Range("D3").Value = "Dr"

# Row 4
Set-TextValue $ws.Range("A4") "05/05/2024"
$ws.Range("B4").Value = "Salary Credit"
Set-TextValue $ws.Range("C4") "50,000.00"
$ws.Range("D4").Value = "Credit"

# Row 5
Set-TextValue $ws.Range("A5") "10-05-2024"
$ws.Range("B5").Value = "Online Shopping"
Set-TextValue $ws.Range("C5") "2,000.00"
$ws.Range("D5").Value = "Debit"
